# Add code 5.0 for Joint Code
$p = $ppt.ActivePresentation

# --- Re-cache the datetimeFigureOut placeholder (9/25/20 -> 9/29/20) on the
# slide master and on every slide layout ---
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
  $sh = $master.Shapes.Item($i)
  if ($sh.Name -like "Date Placeholder*") {
    $sh.TextFrame.TextRange.Text = "9/29/20"
  }
}
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
  $layout = $master.CustomLayouts.Item($i)
  for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
    $sh = $layout.Shapes.Item($j)
    if ($sh.Name -like "Date Placeholder*") {
      $sh.TextFrame.TextRange.Text = "9/29/20"
    }
  }
}

# --- Slide 2 ("Cost Comparison") ---
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Cost Comparison (Double Power Flow)"

$tbl2 = $s2.Shapes.Item(2).Table

# Row 2 ("Real time cost") gains a "(5:00-8:00 Cost)" line, and the three
# stochastic columns gain their own parenthetical delta line.
$tbl2.Cell(2, 1).Shape.TextFrame.TextRange.Text = "Real time cost`r(5:00-8:00 Cost)"
$tbl2.Cell(2, 5).Shape.TextFrame.TextRange.Text = "42,957`r(12,286)"
$tbl2.Cell(2, 6).Shape.TextFrame.TextRange.Text = "43,149`r(12,168)"
$tbl2.Cell(2, 7).Shape.TextFrame.TextRange.Text = "42,778`r(12,179)"

# Row 3 ("Solar Cul Cost") - one value updated
$tbl2.Cell(3, 5).Shape.TextFrame.TextRange.Text = "122.77"

# Row 4 ("Peak Cost") - three values updated
$tbl2.Cell(4, 5).Shape.TextFrame.TextRange.Text = "1,262,500"
$tbl2.Cell(4, 6).Shape.TextFrame.TextRange.Text = "1,638,300"
$tbl2.Cell(4, 7).Shape.TextFrame.TextRange.Text = "1,258,600"

# --- Slide 3 ("Cost Comparison") ---
$s3 = $p.Slides.Item(3)
$s3.SlideShowTransition.Hidden = $true

$tbl3 = $s3.Shapes.Item(2).Table
$tbl3.Cell(4, 5).Shape.TextFrame.TextRange.Text = "12,520,000"
